$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-09 Thursday", "2024-05-10 Friday"),
    @("533÷4=133, 1", "491÷7=70, 1"),
    @("527÷5=105, 2", "190÷9=21, 1"),
    @("843÷2=421, 1", "799÷4=199, 3"),
    @("750÷7=107, 1", "113÷5=22, 3"),
    @("315÷5=63, 0", "165÷4=41, 1"),
    @("677÷9=75, 2", "930÷3=310, 0"),
    @("844÷7=120, 4", "582÷2=291, 0"),
    @("822÷7=117, 3", "530÷2=265, 0"),
    @("949÷6=158, 1", "633÷6=105, 3"),
    @("280÷6=46, 4", "306÷5=61, 1"),
    @("816÷4=204, 0", "826÷9=91, 7"),
    @("161÷7=23, 0", "480÷2=240, 0"),
    @("545÷6=90, 5", "347÷4=86, 3"),
    @("150÷9=16, 6", "271÷3=90, 1"),
    @("314÷8=39, 2", "204÷6=34, 0"),
    @("996÷3=332, 0", "741÷7=105, 6"),
    @("425÷9=47, 2", "654÷6=109, 0"),
    @("630÷2=315, 0", "490÷6=81, 4"),
    @("663÷4=165, 3", "660÷8=82, 4"),
    @("743÷5=148, 3", "604÷3=201, 1"),
    @("776÷3=258, 2", "620÷4=155, 0"),
    @("406÷8=50, 6", "646÷8=80, 6"),
    @("850÷4=212, 2", "123÷4=30, 3"),
    @("607÷3=202, 1", "574÷6=95, 4"),
    @("199÷2=99, 1", "331÷3=110, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
